$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the two runs of the "batch accumulation=32" bullet
# into a single run (same visible text, just no longer split in two).
# ---------------------------------------------------------------------
$mergedText = "loss=cross_entropy AND CropDepth mode=”random” AND batch accumulation=32 AND image depth=48 – equivalent to (5)"

$rng = $d.Content
$rng.Find.Execute($mergedText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeStart = $rng.Start

# Replace via an intermediate placeholder so the "same text" case isn't
# treated as a no-op by the engine.
$rng.Text = "IRON_PLACEHOLDER_0001"
$rng2 = $d.Range($mergeStart, $mergeStart + ("IRON_PLACEHOLDER_0001").Length)
$rng2.Text = $mergedText

# Touch + reset a formatting property on the merged run so it keeps an
# explicit (empty) <w:rPr/> element, matching the canonical output.
$rng3 = $d.Range($mergeStart, $mergeStart + $mergedText.Length)
$rng3.Bold = 1
$rng3.Bold = 0

# ---------------------------------------------------------------------
# Change 2: add new "idea" bullets to the Ideas section.
#   - insert "4.2 Unet with FC ..." right before item "5. ..."
#   - insert items "6.", "8.", "9." right after item "5. ..."
# ---------------------------------------------------------------------
$item5Text = "5. Think of something for active contours in neural networks, conditional random fields, refinement network, boundary loss term"

$item42 = "4.2 Unet with FC from bottleneck to predict classes per voxel. Possibly with 1x1 conv to compress feature maps"
$item6a = "6. Modify elastictransform to center on random organ "
$item6b = "(also alpha / sigma calculation based on organ)"
$item8a = "8. Train longer "
$item8b = "(CE no elastic, FL with beta 0.25)"
$item9  = "9. Cross entropy with class weights based on beta"

# --- insert the "4.2 ..." paragraph right before item 5 ---
$rng = $d.Content
$rng.Find.Execute($item5Text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.InsertParagraphBefore()

$item5Para = $rng.Paragraphs(1)
$newPara1 = $item5Para.Previous()
$newPara1.Range.Text = $item42

# --- insert the "6.", "8.", "9." paragraphs right after item 5 ---
$rng = $d.Content
$rng.Find.Execute($item5Text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.InsertParagraphAfter()

$item5Para = $d.Content
$item5Para.Find.Execute($item5Text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p6 = $item5Para.Paragraphs(1).Next()
$p8 = $p6.Next()
$p9 = $p8.Next()

# paragraph 6: two runs
$p6Start = $p6.Range.Start
$r1 = $d.Range($p6Start, $p6Start)
$r1.InsertAfter($item6a)
$r1End = $r1.End
$r2 = $d.Range($r1End, $r1End)
$r2.InsertAfter($item6b)
$r1Fmt = $d.Range($p6Start, $r1End)
$r1Fmt.Bold = 1
$r1Fmt.Bold = 0

# paragraph 8: two runs
$p8Start = $p8.Range.Start
$r1 = $d.Range($p8Start, $p8Start)
$r1.InsertAfter($item8a)
$r1End = $r1.End
$r2 = $d.Range($r1End, $r1End)
$r2.InsertAfter($item8b)
$r1Fmt = $d.Range($p8Start, $r1End)
$r1Fmt.Bold = 1
$r1Fmt.Bold = 0

# paragraph 9: single run
$p9.Range.Text = $item9
